$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header columns: Email, Games, Name, Age
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Games"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Age"

# Move the Games dropdown list validation from column D to column B
$ws.Range("D2:D100000").Validation.Delete()
$ws.Range("B2:B100000").Validation.Add(3, 1, 1, '"Super Mario,SONIC,Zelda,GTA"')
$ws.Range("B2:B100000").Validation.IgnoreBlank = $true
